$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 86989.5
$ws.Range("J95").Value = 86989.5
$ws.Range("L95").Value = 86989.5
$ws.Range("N95").Value = -92481.5
$ws.Range("H127").Value = 2993.9546
$ws.Range("I127").Value = 3145.85
$ws.Range("J127").Value = 1475
$ws.Range("K127").Value = 9437.549999999999
$ws.Range("L127").Value = 4425
$ws.Range("M127").Value = -4477.549999999999
$ws.Range("N127").Value = -14345
$ws.Range("H131").Value = 1373.8667
$ws.Range("I131").Value = 816
$ws.Range("K131").Value = 2448
$ws.Range("M131").Value = 2592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 761.6667
$ws.Range("I2").Value = 768
$ws.Range("J2").Value = 730
$ws.Range("K2").Value = 768
$ws.Range("L2").Value = 730
$ws.Range("M2").Value = -655
$ws.Range("N2").Value = -956
$ws.Range("H61").Value = 3722.625
$ws.Range("I61").Value = 3798.72
$ws.Range("J61").Value = 3450.8572
$ws.Range("K61").Value = 3798.72
$ws.Range("L61").Value = 3450.8572
$ws.Range("M61").Value = -3586.72
$ws.Range("N61").Value = -3874.8572
$ws.Range("H74").Value = 11859.75
$ws.Range("I74").Value = 1823.875
$ws.Range("K74").Value = 1823.875
$ws.Range("M74").Value = -949.875
$ws.Range("H77").Value = 11859.75
$ws.Range("I77").Value = 1823.875
$ws.Range("K77").Value = 9119.375
$ws.Range("M77").Value = -4751.375
$ws.Range("H97").Value = 2045.5555
$ws.Range("I97").Value = 1328
$ws.Range("K97").Value = 1328
$ws.Range("M97").Value = -832
$ws.Range("H102").Value = 4318.231
$ws.Range("I102").Value = 4398.2085
$ws.Range("K102").Value = 4398.2085
$ws.Range("M102").Value = -2776.2085
$ws.Range("H116").Value = 761.6667
$ws.Range("I116").Value = 768
$ws.Range("J116").Value = 730
$ws.Range("K116").Value = 768
$ws.Range("L116").Value = 730
$ws.Range("M116").Value = 1526
$ws.Range("N116").Value = -5318
$ws.Range("H122").Value = 2111.5789
$ws.Range("I122").Value = 1971
$ws.Range("J122").Value = 2416.1667
$ws.Range("K122").Value = 5913
$ws.Range("L122").Value = 7248.500100000001
$ws.Range("M122").Value = -3463
$ws.Range("N122").Value = -12148.5001
$ws.Range("H132").Value = 2274528
$ws.Range("I132").Value = 3572977.2
$ws.Range("J132").Value = 2241.75
$ws.Range("K132").Value = 10718931.6
$ws.Range("L132").Value = 6725.25
$ws.Range("M132").Value = -10716401.6
$ws.Range("N132").Value = -11785.25
$ws.Range("H133").Value = 76326.664
$ws.Range("J133").Value = 76326.664
$ws.Range("L133").Value = 76326.664
$ws.Range("N133").Value = -81386.664
$ws.Range("H134").Value = 66998.60000000001
$ws.Range("J134").Value = 66998.60000000001
$ws.Range("L134").Value = 66998.60000000001
$ws.Range("N134").Value = -77138.60000000001
$ws.Range("H136").Value = 3722.625
$ws.Range("I136").Value = 3798.72
$ws.Range("J136").Value = 3450.8572
$ws.Range("K136").Value = 11396.16
$ws.Range("L136").Value = 10352.5716
$ws.Range("M136").Value = -8846.16
$ws.Range("N136").Value = -15452.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 761.6667
$ws.Range("I3").Value = 768
$ws.Range("J3").Value = 730
$ws.Range("K3").Value = 768
$ws.Range("L3").Value = 730
$ws.Range("M3").Value = -654
$ws.Range("N3").Value = -958
$ws.Range("H86").Value = 1965.9412
$ws.Range("I86").Value = 1796.7273
$ws.Range("J86").Value = 2276.1667
$ws.Range("K86").Value = 1796.7273
$ws.Range("L86").Value = 2276.1667
$ws.Range("M86").Value = -673.7273
$ws.Range("N86").Value = -4522.1667
$ws.Range("H89").Value = 1965.9412
$ws.Range("I89").Value = 1796.7273
$ws.Range("J89").Value = 2276.1667
$ws.Range("K89").Value = 8983.636500000001
$ws.Range("L89").Value = 11380.8335
$ws.Range("M89").Value = -3367.636500000001
$ws.Range("N89").Value = -22612.8335
$ws.Range("H94").Value = 2804.8
$ws.Range("I94").Value = 2424.44
$ws.Range("K94").Value = 2424.44
$ws.Range("M94").Value = -1973.44
$ws.Range("H103").Value = 43833.332
$ws.Range("J103").Value = 43833.332
$ws.Range("L103").Value = 43833.332
$ws.Range("N103").Value = -46177.332
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H134").Value = 30253
$ws.Range("I134").Value = 8006
$ws.Range("J134").Value = 52500
$ws.Range("K134").Value = 24018
$ws.Range("L134").Value = 157500
$ws.Range("M134").Value = -21483
$ws.Range("N134").Value = -162570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 4856.8335
$ws.Range("J31").Value = 3081.2856
$ws.Range("K31").Value = 4856.8335
$ws.Range("L31").Value = 3081.2856
$ws.Range("M31").Value = -4561.8335
$ws.Range("N31").Value = -3671.2856
$ws.Range("I34").Value = 4856.8335
$ws.Range("J34").Value = 3081.2856
$ws.Range("K34").Value = 4856.8335
$ws.Range("L34").Value = 3081.2856
$ws.Range("M34").Value = -4654.8335
$ws.Range("N34").Value = -3485.2856
$ws.Range("H134").Value = 3445.6924
$ws.Range("I134").Value = 2459.5
$ws.Range("K134").Value = 7378.5
$ws.Range("M134").Value = -4843.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2085634.5
$ws.Range("J121").Value = 3574941
$ws.Range("L121").Value = 10724823
$ws.Range("N121").Value = -10727443
$ws.Range("H131").Value = 1998.8
$ws.Range("J131").Value = 1998.8
$ws.Range("L131").Value = 5996.4
$ws.Range("N131").Value = -16076.4
$ws.Range("H139").Value = 5636.3
$ws.Range("I139").Value = 4623.2856
$ws.Range("K139").Value = 13869.8568
$ws.Range("M139").Value = -8729.856800000001
$ws.Range("H140").Value = 1228.125
$ws.Range("I140").Value = 1228.125
$ws.Range("K140").Value = 3684.375
$ws.Range("M140").Value = 1495.625
$ws.Range("H141").Value = 11319.071
$ws.Range("I141").Value = 10315.182
$ws.Range("K141").Value = 30945.546
$ws.Range("M141").Value = -25765.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2660.6
$ws.Range("I5").Value = 3275.75
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 3275.75
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -3163.75
$ws.Range("N5").Value = -424
$ws.Range("H25").Value = 2375
$ws.Range("J25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("N25").Value = -3558
$ws.Range("H122").Value = 3146.5
$ws.Range("I122").Value = 2556.5
$ws.Range("K122").Value = 7669.5
$ws.Range("M122").Value = -5219.5
$ws.Range("H132").Value = 17061.566
$ws.Range("I132").Value = 18074.893
$ws.Range("K132").Value = 54224.679
$ws.Range("M132").Value = -51694.679

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3305.7917
$ws.Range("I40").Value = 3273.5715
$ws.Range("J40").Value = 3531.3333
$ws.Range("K40").Value = 3273.5715
$ws.Range("L40").Value = 3531.3333
$ws.Range("M40").Value = -3137.5715
$ws.Range("N40").Value = -3803.3333
$ws.Range("H106").Value = 19500
$ws.Range("J106").Value = 19500
$ws.Range("L106").Value = 19500
$ws.Range("N106").Value = -22024
$ws.Range("H122").Value = 3920.077
$ws.Range("I122").Value = 4415.5
$ws.Range("K122").Value = 13246.5
$ws.Range("M122").Value = -10796.5
$ws.Range("H132").Value = 4169758
$ws.Range("I132").Value = 10002198
$ws.Range("K132").Value = 30006594
$ws.Range("M132").Value = -30004064
$ws.Range("H136").Value = 11192.538
$ws.Range("I136").Value = 4734.5
$ws.Range("K136").Value = 14203.5
$ws.Range("M136").Value = -11653.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1537.6
$ws.Range("I132").Value = 1443.5834
$ws.Range("K132").Value = 4330.7502
$ws.Range("M132").Value = -1800.7502
$ws.Range("H136").Value = 388892.47
$ws.Range("I136").Value = 533328.4
$ws.Range("J136").Value = 3730
$ws.Range("K136").Value = 1599985.2
$ws.Range("L136").Value = 11190
$ws.Range("M136").Value = -1597435.2
$ws.Range("N136").Value = -16290
